$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at 68-69 (weekly update: new Primera/Segunda record),
# pushing the former rows 68-85 down to 70-87.
$ws.Range("A68:A69").EntireRow.Insert()

# Row 68: new "Primera" quality record for 2023-08-16 (serial 45154)
$ws.Range("A68").Value = 3
$ws.Range("B68").Value = "Femacal de La Calera"
$ws.Range("C68").Value = "Coquimbo"
$ws.Range("D68").Value = 45154
$ws.Range("D68").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E68").Value = 5
$ws.Range("F68").Value = "Fruta"
$ws.Range("G68").Value = 100108
$ws.Range("H68").Value = "Tropicales y subtropicales"
$ws.Range("I68").Value = 100108004
$ws.Range("J68").Value = "Papaya"
$ws.Range("K68").Value = "Cultivar IV Región"
$ws.Range("L68").Value = "Primera"
$ws.Range("M68").Value = 54
$ws.Range("N68").Value = 20000
$ws.Range("O68").Value = 20000
$ws.Range("P68").Value = 20000
$ws.Range("Q68").Value = "$/bandeja 10 kilos"
$ws.Range("R68").Value = "Provincia del Elquí"
$ws.Range("S68").Value = 2000
$ws.Range("T68").Value = 10

# Row 69: new "Segunda" quality record for 2023-08-16 (serial 45154)
$ws.Range("A69").Value = 3
$ws.Range("B69").Value = "Femacal de La Calera"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 45154
$ws.Range("D69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E69").Value = 5
$ws.Range("F69").Value = "Fruta"
$ws.Range("G69").Value = 100108
$ws.Range("H69").Value = "Tropicales y subtropicales"
$ws.Range("I69").Value = 100108004
$ws.Range("J69").Value = "Papaya"
$ws.Range("K69").Value = "Cultivar IV Región"
$ws.Range("L69").Value = "Segunda"
$ws.Range("M69").Value = 56
$ws.Range("N69").Value = 17000
$ws.Range("O69").Value = 17000
$ws.Range("P69").Value = 17000
$ws.Range("Q69").Value = "$/bandeja 10 kilos"
$ws.Range("R69").Value = "Provincia del Elquí"
$ws.Range("S69").Value = 1700
$ws.Range("T69").Value = 10
